$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (and the Solana/WrappedEther
# row-order swap) as refreshed by the scheduled GitHub Actions scrape.
# Values are written with a leading apostrophe so Excel keeps them as text
# (these columns hold formatted/locale-style numeric text, not numbers),
# then ClearFormats() strips the implicit "quote prefix" text format Excel
# applies so the cell style matches the rest of the unstyled data cells.
$updates = @(
    @{ Cell = "D2"; Value = "30.471.47" },
    @{ Cell = "E2"; Value = "  -0.85%  " },
    @{ Cell = "D3"; Value = "1.894.61" },
    @{ Cell = "E4"; Value = "  -0.18%  " },
    @{ Cell = "D5"; Value = "237.98" },
    @{ Cell = "E5"; Value = "  +0.72%  " },
    @{ Cell = "D6"; Value = "1.000" },
    @{ Cell = "E6"; Value = "  -0.19%  " },
    @{ Cell = "D7"; Value = "0.4894" },
    @{ Cell = "E7"; Value = "  +0.27%  " },
    @{ Cell = "D8"; Value = "0.2934" },
    @{ Cell = "E8"; Value = "  +1.19%  " },
    @{ Cell = "D9"; Value = "0.06686" },
    @{ Cell = "B10"; Value = "WrappedEther" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" },
    @{ Cell = "D10"; Value = "1.881.29" },
    @{ Cell = "E10"; Value = "  -0.51%  " },
    @{ Cell = "B11"; Value = "Solana" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol" },
    @{ Cell = "D11"; Value = "17.09" },
    @{ Cell = "E11"; Value = "  +2.31%  " },
    @{ Cell = "E12"; Value = "  +1.35%  " },
    @{ Cell = "D13"; Value = "5.142" },
    @{ Cell = "E13"; Value = "  +2.58%  " },
    @{ Cell = "D14"; Value = "88.11" },
    @{ Cell = "E14"; Value = "  -1.35%  " },
    @{ Cell = "D15"; Value = "0.6658" },
    @{ Cell = "E15"; Value = "  +0.21%  " },
    @{ Cell = "D16"; Value = "30.443.60" },
    @{ Cell = "E16"; Value = "  -0.78%  " },
    @{ Cell = "D17"; Value = "13.48" },
    @{ Cell = "E17"; Value = "  +3.84%  " },
    @{ Cell = "D18"; Value = "0.000007827" },
    @{ Cell = "E18"; Value = "  -0.74%  " },
    @{ Cell = "E19"; Value = "  -0.13%  " },
    @{ Cell = "D20"; Value = "2.141.31" },
    @{ Cell = "E20"; Value = "  +0.51%  " },
    @{ Cell = "D21"; Value = "5.319" },
    @{ Cell = "E21"; Value = "  +12.15%  " },
    @{ Cell = "D22"; Value = "1.001" },
    @{ Cell = "E22"; Value = "  -0.17%  " },
    @{ Cell = "D23"; Value = "189.62" },
    @{ Cell = "E23"; Value = "  -1.06%  " },
    @{ Cell = "D24"; Value = "6.152" },
    @{ Cell = "E24"; Value = "  +1.44%  " },
    @{ Cell = "D25"; Value = "9.490" },
    @{ Cell = "E25"; Value = "  +2.08%  " },
    @{ Cell = "D26"; Value = "164.86" },
    @{ Cell = "E26"; Value = "  +3.37%  " },
    @{ Cell = "D27"; Value = "18.31" },
    @{ Cell = "E27"; Value = "  -0.02%  " },
    @{ Cell = "D28"; Value = "1.935" },
    @{ Cell = "E28"; Value = "  +6.13%  " },
    @{ Cell = "D29"; Value = "1.463" },
    @{ Cell = "E29"; Value = "  +4.58%  " },
    @{ Cell = "D30"; Value = "4.359" },
    @{ Cell = "E30"; Value = "  +2.16%  " },
    @{ Cell = "D31"; Value = "0.09175" },
    @{ Cell = "E31"; Value = "  +1.58%  " },
    @{ Cell = "D32"; Value = "4.084" },
    @{ Cell = "E32"; Value = "  +3.61%  " },
    @{ Cell = "D33"; Value = "0.05216" },
    @{ Cell = "E33"; Value = "  +0.15%  " },
    @{ Cell = "D34"; Value = "0.7427" },
    @{ Cell = "E34"; Value = "  +1.41%  " },
    @{ Cell = "D35"; Value = "1.100" },
    @{ Cell = "E35"; Value = "  +1.31%  " },
    @{ Cell = "D36"; Value = "2.716" },
    @{ Cell = "E36"; Value = "  +0.75%  " },
    @{ Cell = "D37"; Value = "0.01819" },
    @{ Cell = "E37"; Value = "  -0.16%  " },
    @{ Cell = "D38"; Value = "2.674" },
    @{ Cell = "E38"; Value = "  +0.06%  " },
    @{ Cell = "D39"; Value = "0.9186" },
    @{ Cell = "E39"; Value = "  -0.50%  " },
    @{ Cell = "D40"; Value = "2.046" },
    @{ Cell = "E40"; Value = "  -0.71%  " },
    @{ Cell = "D41"; Value = "0.4403" },
    @{ Cell = "E41"; Value = "  -1.03%  " },
    @{ Cell = "D42"; Value = "5.949" },
    @{ Cell = "E42"; Value = "  +3.92%  " },
    @{ Cell = "D43"; Value = "106.16" },
    @{ Cell = "E43"; Value = "  +1.42%  " },
    @{ Cell = "D44"; Value = "0.9933" },
    @{ Cell = "E44"; Value = "  -0.76%  " },
    @{ Cell = "D45"; Value = "0.1386" },
    @{ Cell = "E45"; Value = "  +3.88%  " },
    @{ Cell = "D46"; Value = "68.22" },
    @{ Cell = "E46"; Value = "  +18.98%  " },
    @{ Cell = "D47"; Value = "7.615" },
    @{ Cell = "E47"; Value = "  +4.18%  " },
    @{ Cell = "D48"; Value = "9.008" },
    @{ Cell = "E48"; Value = "  +4.16%  " },
    @{ Cell = "E49"; Value = "  +5.06%  " },
    @{ Cell = "D50"; Value = "0.05825" },
    @{ Cell = "E50"; Value = "  -0.20%  " },
    @{ Cell = "D51"; Value = "0.3953" },
    @{ Cell = "E51"; Value = "  -6.46%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.ClearFormats()
}
